# 55Cnc_e workbook update — "Finalizando gráficos para plages"
# Adds a new interpolation ("Temp") table in columns I:P (rows 3-8),
# a small extracted row (P12:T12), a few scratch values (L/N:R 15-19),
# and finishes the page setup for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New header cell: J1 = 5172 (matches the B1:F1 wavelength header row)
# ---------------------------------------------------------------------
$ws.Range("J1").Value = 5172

# ---------------------------------------------------------------------
# 2) Pre-format the new "Temp" block (L3:Y8) with the same style used by
#    the rest of the data table (A2 carries the "Normal_Sheet1" style),
#    so blank trailing cells Q:Y also end up carrying that style.
# ---------------------------------------------------------------------
$ws.Range("L3:Y8").Style = $ws.Range("A2").Style

# Row 3: label + first data row (plain values, no interpolation formulas)
$ws.Range("K3").Value = "Temp"
$ws.Range("L3").Value = 452.1
$ws.Range("M3").Value = 521.5
$ws.Range("N3").Value = 563.5
$ws.Range("O3").Value = 614
$ws.Range("P3").Value = 664.6

# Row 4: first interpolation row (distinct, non-shared formulas)
$ws.Range("I4").Formula = "=J4-100"
$ws.Range("J4").Formula = "=K4*100/5172"
$ws.Range("K4").Value = 5272
$ws.Range("L4").Value = 375.69017542393726
$ws.Range("M4").Value = 363.35132060638654
$ws.Range("N4").Value = 357.75265758453133
$ws.Range("O4").Value = 352.45475500789161
$ws.Range("P4").Value = 347.19488560286391

# Rows 5-8: remaining interpolation rows, filled down as shared formulas
$ws.Range("I5:I8").Formula = "=J5-100"
$ws.Range("J5:J8").Formula = "=K5*100/5172"

$ws.Range("K5").Value = 5472
$ws.Range("L5").Value = 368.90724514904161
$ws.Range("M5").Value = 357.84500566993228
$ws.Range("N5").Value = 352.79719981173849
$ws.Range("O5").Value = 348.03170535124826
$ws.Range("P5").Value = 343.21001157433528

$ws.Range("K6").Value = 5672
$ws.Range("L6").Value = 361.39575297344794
$ws.Range("M6").Value = 351.87075718656757
$ws.Range("N6").Value = 347.47377683108647
$ws.Range("O6").Value = 343.32735153219664
$ws.Range("P6").Value = 339.00668909880858

$ws.Range("K7").Value = 6172
$ws.Range("L7").Value = 339.85551396442656
$ws.Range("M7").Value = 335.19102897627829
$ws.Range("N7").Value = 332.80849857331151
$ws.Range("O7").Value = 330.54203562676855
$ws.Range("P7").Value = 327.71229081007823

$ws.Range("K8").Value = 7172
$ws.Range("L8").Value = 289.6882987073468
$ws.Range("M8").Value = 297.23763878608531
$ws.Range("N8").Value = 299.90134729795949
$ws.Range("O8").Value = 302.28481398220094
$ws.Range("P8").Value = 303.08320393501731

# ---------------------------------------------------------------------
# 3) Extracted row for the chart (keeps the styled "Normal_Sheet1" look)
# ---------------------------------------------------------------------
$ws.Range("P12:T12").Style = $ws.Range("A2").Style
$ws.Range("P12").Value = 347.19488560286391
$ws.Range("Q12").Value = 343.21001157433528
$ws.Range("R12").Value = 339.00668909880858
$ws.Range("S12").Value = 327.71229081007823
$ws.Range("T12").Value = 303.08320393501731

# ---------------------------------------------------------------------
# 4) Scratch values (plain, unstyled) pasted below the table
# ---------------------------------------------------------------------
$ws.Range("L15").Value = 1.9334880123743261
$ws.Range("N15").Value = 1.9334880123743261
$ws.Range("O15").Value = 5.8004640371229641
$ws.Range("P15").Value = 9.6674400618716163
$ws.Range("Q15").Value = 19.334880123743233
$ws.Range("R15").Value = 38.669760247486465

$ws.Range("L16").Value = 5.8004640371229641
$ws.Range("L17").Value = 9.6674400618716163
$ws.Range("L18").Value = 19.334880123743233
$ws.Range("L19").Value = 38.669760247486465

# ---------------------------------------------------------------------
# 5) Selection follows the newly pasted chart-source range
# ---------------------------------------------------------------------
$ws.Range("P12:T12").Select()

# ---------------------------------------------------------------------
# 6) Finish page setup for printing
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
